$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27, shifting existing rows 27-57 down to 28-58.
$ws.Rows("27:27").Insert()

# Populate the newly inserted row 27 with the new record.
$ws.Range("A27").Value = 3
$ws.Range("B27").Value = "Femacal de La Calera"
$ws.Range("C27").Value = "Coquimbo"
$ws.Range("D27").Value = 44902
$ws.Range("E27").Value = 5
$ws.Range("F27").Value = 300000000
$ws.Range("G27").Value = "Espárragos"
$ws.Range("H27").Value = "Verde"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 1250
$ws.Range("K27").Value = 1400
$ws.Range("L27").Value = 1500
$ws.Range("M27").Value = 1452
$ws.Range("N27").Value = "$/kilo"
$ws.Range("O27").Value = "Provincia de Quillota"
$ws.Range("P27").Value = 1452
$ws.Range("Q27").Value = 1
$ws.Range("R27").Value = "Hortaliza"
